# hla-drb1 excel: add a "variant_effect" tab summarising a single VEP
# (Variant Effect Predictor) annotation row, and leave the "vcf color"
# sheet's selection where the author last clicked (D21) instead of its
# old B105:D109 block selection.

$wb = $excel.ActiveWorkbook

# --- Leave a trail on "vcf color": author had simply clicked on D21
# before switching away to review the new tab. -----------------------
$vcfColor = $wb.Worksheets.Item("vcf color")
$vcfColor.Activate()
$vcfColor.Range("D21").Select()

# --- Add the new "variant_effect" worksheet after the existing tabs. -
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "variant_effect"

# Single VEP annotation row for the stop-gained call at 12:2157346.
$ws.Range("A1").Value = "12_2157346_G/*/T"
$ws.Range("B1").Value = "12:2157346"
$ws.Range("C1").Value = "T"
$ws.Range("D1").Value = 474860
$ws.Range("E1").Value = "NM_001014768.1"
$ws.Range("F1").Value = "Transcript"
$ws.Range("G1").Value = "stop_gained"
$ws.Range("H1").Value = 152
$ws.Range("I1").Value = 114
$ws.Range("J1").Value = 38
$ws.Range("K1").Value = "Y/*"
$ws.Range("L1").Value = "taC/taA"
$ws.Range("M1").Value = "-"
$ws.Range("N1").Value = "IMPACT=HIGH;STRAND=-1"

# Touching NumberFormat materialises an explicit (Normal-equivalent)
# style record for the row, matching the extra cellXfs entry that shows
# up in the saved workbook.
$ws.Range("A1:N1").NumberFormat = "General"

# Make the new tab the active one, with the full header row selected,
# mirroring the author's final view before saving.
$ws.Activate()
$ws.Range("A1:N1").Select()
